# Updated cryptos list on Sat Jun  8 10:45:29 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $origStyle = $r.Style
    # Force text interpretation so numeric-looking strings (e.g. "687.59",
    # "1.00", "69.378.74") are kept as exact text instead of being
    # parsed/rounded into floating point numbers.
    $r.NumberFormat = "@"
    $r.Value = $value
    # Restore the original (unset/"Normal") cell style so no stray
    # direct formatting is introduced.
    $r.Style = $origStyle
}

# Row 2 - Bitcoin
Set-TextValue "D2" "69.378.74"
Set-TextValue "E2" "  -2.61%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.683.71"
Set-TextValue "E3" "  -3.19%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.07%  "

# Row 5 - BNB
Set-TextValue "D5" "687.59"
Set-TextValue "E5" "  -2.07%  "

# Row 6 - Solana
Set-TextValue "D6" "161.78"
Set-TextValue "E6" "  -5.44%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.680.90"
Set-TextValue "E7" "  -3.26%  "

# Row 9 - XRP
Set-TextValue "E9" "  -5.49%  "

# Row 10 - Dogecoin
Set-TextValue "E10" "  -8.34%  "

# Row 11 - Toncoin
Set-TextValue "D11" "7.37"
Set-TextValue "E11" "  -1.70%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.439"
Set-TextValue "E12" "  -7.52%  "

# Row 13 - ShibaInu
Set-TextValue "E13" "  -5.70%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "4.307.07"
Set-TextValue "E14" "  -3.08%  "

# Row 15 - Avalanche
Set-TextValue "D15" "33.08"

# Row 16 - WrappedEther
Set-TextValue "D16" "3.690.54"
Set-TextValue "E16" "  -4.28%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "69.407.95"
Set-TextValue "E17" "  -2.52%  "

# Row 18 - TRON
Set-TextValue "E18" "  -1.53%  "

# Row 19 - Chainlink
Set-TextValue "E19" "  -7.89%  "

# Row 20 - Polkadot
Set-TextValue "E20" "  -8.80%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "476.30"
Set-TextValue "E21" "  -7.32%  "

# Row 22 - Uniswap
Set-TextValue "D22" "9.95"
Set-TextValue "E22" "  -5.35%  "

# Row 23 - Polygon
Set-TextValue "E23" "  -7.60%  "

# Row 24 - Litecoin
Set-TextValue "E24" "  -4.91%  "

# Row 25 - WrappedeETH
Set-TextValue "D25" "3.829.10"
Set-TextValue "E25" "  -3.16%  "

# Row 26 - PEPE
Set-TextValue "E26" "  -9.41%  "

# Row 27 - Dai
Set-TextValue "E27" "  +0.07%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue "D28" "11.25"
Set-TextValue "E28" "  -7.11%  "

# Row 29 - RenderToken
Set-TextValue "D29" "9.34"
Set-TextValue "E29" "  -9.54%  "

# Row 30 - Fetch.AI
Set-TextValue "E30" "  -10.63%  "

# Row 31 - PancakeSwap
Set-TextValue "E31" "  -10.13%  "

# Row 32 - NEARProtocol
Set-TextValue "E32" "  -7.51%  "

# Row 33 - ImmutableX
Set-TextValue "E33" "  -7.79%  "

# Row 34 - was Binance-PegBSC-USD, now Kaspa (rows 34/35 content swapped)
Set-TextValue "B34" "Kaspa"
Set-TextValue "C34" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D34" "0.167"
Set-TextValue "E34" "  -4.21%  "

# Row 35 - was Kaspa, now Binance-PegBSC-USD
Set-TextValue "B35" "Binance-PegBSC-USD"
Set-TextValue "C35" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D35" "1.00"
Set-TextValue "E35" "  +0.05%  "

# Row 36 - EthereumClassic
Set-TextValue "D36" "26.87"
Set-TextValue "E36" "  -7.41%  "

# Row 37 - RenzoRestakedETH
Set-TextValue "D37" "3.650.75"
Set-TextValue "E37" "  -3.06%  "

# Row 38 - Aptos
Set-TextValue "E38" "  -8.09%  "

# Row 39 - Filecoin
Set-TextValue "E39" "  -1.52%  "

# Row 40 - Stacks
Set-TextValue "E40" "  -2.52%  "

# Row 41 - Hedera
Set-TextValue "E41" "  -8.75%  "

# Row 43 - FirstDigitalUSD
Set-TextValue "E43" "  +0.00%  "

# Row 44 - Mantle
Set-TextValue "E44" "  -6.42%  "

# Row 45 - Monero
Set-TextValue "D45" "163.67"
Set-TextValue "E45" "  -5.48%  "

# Row 46 - OKB
Set-TextValue "D46" "48.11"
Set-TextValue "E46" "  -2.71%  "

# Row 47 - InjectiveProtocol
Set-TextValue "D47" "29.77"
Set-TextValue "E47" "  +1.23%  "

# Row 48 - was SuiNetwork, now dogwifhat (rows 48/49 content swapped)
Set-TextValue "B48" "dogwifhat"
Set-TextValue "C48" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D48" "2.75"
Set-TextValue "E48" "  -15.50%  "

# Row 49 - was dogwifhat, now SuiNetwork
Set-TextValue "B49" "SuiNetwork"
Set-TextValue "C49" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D49" "1.13"
Set-TextValue "E49" "  -1.32%  "

# Row 50 - FLOKI
Set-TextValue "D50" "0.000282"
Set-TextValue "E50" "  -8.35%  "

# Row 51 - ONDO
Set-TextValue "E51" "  -4.37%  "
